# Update the "Fecha" column (A2:A23) from 2025-05-15 (45792) to 2025-06-17 (45825)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:A23").Value = 45825
